$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.99828577547460617
$ws.Range("BO1").Value = 0.754174144953905
$ws.Range("BP1").Value = 0.70511613568263609
$ws.Range("A2").Value = 0.70303173535771024
$ws.Range("C2").Value = 0.72844068033737464
$ws.Range("D2").Value = 0.57231144135010115
$ws.Range("BP2").Value = 0.90976282809818232
$ws.Range("C4").Value = 0.80654828439374404
$ws.Range("E4").Value = 0.86887771086717724
$ws.Range("AU4").Value = 0.86821051467633303
$ws.Range("C5").Value = 0.57360488653202291
$ws.Range("F5").Value = 0.89596810102910651
$ws.Range("AO5").Value = 0.91347554484930538
$ws.Range("D6").Value = 0.78269689057981051
$ws.Range("F7").Value = 0.70212629207563948
$ws.Range("AK7").Value = 0.9863406027914996
$ws.Range("E8").Value = 0.81801494810644204
$ws.Range("F8").Value = 0.85971721867524531
$ws.Range("G8").Value = 0.77012186923877324
$ws.Range("J8").Value = 0.7726280348804746
$ws.Range("K9").Value = 0.73594594102690292
$ws.Range("AW9").Value = 0.72180620700581777
$ws.Range("K10").Value = 0.86089167528985211
$ws.Range("M11").Value = 0.62688876337323673
$ws.Range("M12").Value = 0.82298593288775646
$ws.Range("N12").Value = 0.9719502886205057
$ws.Range("W13").Value = 0.75853277876851677
$ws.Range("M14").Value = 0.64545993266656421
$ws.Range("N15").Value = 0.91900756099613878
$ws.Range("P15").Value = 0.98630314260204188
$ws.Range("Q15").Value = 0.69954808414180603
$ws.Range("N16").Value = 0.99582967142846146
$ws.Range("Q16").Value = 0.8313203494363075
$ws.Range("Q18").Value = 0.64015412252493786
$ws.Range("R19").Value = 0.74577242429555435
$ws.Range("U19").Value = 0.69322930884961242
$ws.Range("AD19").Value = 0.8580618738488307
$ws.Range("R20").Value = 0.72738144273843641
$ws.Range("X20").Value = 0.97581572742787037
$ws.Range("W21").Value = 0.82924603998019664
$ws.Range("AD21").Value = 0.7857185407965892
$ws.Range("P22").Value = 0.93187854591150998
$ws.Range("T22").Value = 0.85175604087254575
$ws.Range("AB22").Value = 0.93951848561134788
$ws.Range("V23").Value = 0.83843010686868169
$ws.Range("V24").Value = 0.96202227346672764
$ws.Range("W25").Value = 0.96951374399840051
$ws.Range("X25").Value = 0.81591252633256195
$ws.Range("AA25").Value = 0.9287102066172972
$ws.Range("D26").Value = 0.89856281790095416
$ws.Range("K26").Value = 0.90998311145763511
$ws.Range("X26").Value = 0.80668440390423823
$ws.Range("Y26").Value = 0.86120182712697613
$ws.Range("AB26").Value = 0.71002137405750188
$ws.Range("J27").Value = 0.87701589081309406
$ws.Range("AB27").Value = 0.94518876854336642
$ws.Range("BC27").Value = 0.79935801427401487
$ws.Range("M28").Value = 0.70067009871941344
$ws.Range("AB29").Value = 0.96153903967709486
$ws.Range("D30").Value = 0.64925163722237977
$ws.Range("AC30").Value = 0.68165735427150209
$ws.Range("W31").Value = 0.89182257807799026
$ws.Range("AC31").Value = 0.91460764404555628
$ws.Range("AS31").Value = 0.89630410726717291
$ws.Range("AG32").Value = 0.91216300708269227
$ws.Range("AH32").Value = 0.80782072506324654
$ws.Range("AE33").Value = 0.84623341663100038
$ws.Range("AI33").Value = 0.96377621370789301
$ws.Range("AR33").Value = 0.52268856968981436
$ws.Range("BJ33").Value = 0.59653101165679245
$ws.Range("AI34").Value = 0.95851614295975041
$ws.Range("AJ34").Value = 0.91123222877881727
$ws.Range("AJ35").Value = 0.81771391530359328
$ws.Range("G36").Value = 0.97969424591308607
$ws.Range("AL36").Value = 0.95433044155043278
$ws.Range("AI37").Value = 0.84451259052747085
$ws.Range("AK38").Value = 0.80368985519704128
$ws.Range("AP38").Value = 0.97449440139123777
$ws.Range("BI39").Value = 0.90648231460065332
$ws.Range("AK40").Value = 0.64314810793110477
$ws.Range("AM41").Value = 0.94125277351428038
$ws.Range("AN41").Value = 0.84501997295949116
$ws.Range("AQ41").Value = 0.92121335163667217
$ws.Range("AN42").Value = 0.85860249796022536
$ws.Range("AO42").Value = 0.99981109238863963
$ws.Range("AR42").Value = 0.79213579092962738
$ws.Range("AP43").Value = 0.86282171591465429
$ws.Range("AR43").Value = 0.87997355866695703
$ws.Range("AS43").Value = 0.63089420513672101
$ws.Range("AR45").Value = 0.87571687161571199
$ws.Range("AR46").Value = 0.89343241891399972
$ws.Range("AG47").Value = 0.5759076346276748
$ws.Range("AS47").Value = 0.86659165252714154
$ws.Range("AT47").Value = 0.97207178328091504
$ws.Range("AT48").Value = 0.86951550470470962
$ws.Range("AV49").Value = 0.96312791691848643
$ws.Range("AY49").Value = 0.84873289309200406
$ws.Range("W50").Value = 0.9250502052350158
$ws.Range("AV50").Value = 0.6986728605282575
$ws.Range("AW50").Value = 0.80086314112351631
$ws.Range("D51").Value = 0.69335414966220354
$ws.Range("AX51").Value = 0.71800743812401324
$ws.Range("BA51").Value = 0.9060243448670392
$ws.Range("D52").Value = 0.71534472344688393
$ws.Range("BA52").Value = 0.91684808795884831
$ws.Range("BC53").Value = 0.5161182476750672
$ws.Range("B54").Value = 0.91687916922947577
$ws.Range("BA54").Value = 0.73933604737827729
$ws.Range("T55").Value = 0.90486696222199225
$ws.Range("BB55").Value = 0.79658314614446346
$ws.Range("BD55").Value = 0.84194436022268104
$ws.Range("W56").Value = 0.75790028093167006
$ws.Range("BF57").Value = 0.96615090670605985
$ws.Range("BG57").Value = 0.94963001356646037
$ws.Range("AL58").Value = 0.79992673029475969
$ws.Range("BD58").Value = 0.70173241193108926
$ws.Range("BH58").Value = 0.91162182235264089
$ws.Range("BM58").Value = 0.91377039901967738
$ws.Range("BG60").Value = 0.73742483875305309
$ws.Range("BI60").Value = 0.7657907615169488
$ws.Range("BJ60").Value = 0.79397069741565474
$ws.Range("Q61").Value = 0.98933542128628682
$ws.Range("AH61").Value = 0.9654419043166822
$ws.Range("BG61").Value = 0.96821210878717068
$ws.Range("BK61").Value = 0.95305120059313164
$ws.Range("BK62").Value = 0.83697652285103152
$ws.Range("BI64").Value = 0.68693692487451419
$ws.Range("BK64").Value = 0.56990948921074158
$ws.Range("BM64").Value = 0.59993739481043995
$ws.Range("BN64").Value = 0.99323285610133016
$ws.Range("F65").Value = 0.98860212738327913
$ws.Range("S65").Value = 0.95287165322746459
$ws.Range("BK65").Value = 0.82484506774366229
$ws.Range("BN65").Value = 0.80006302341358904
$ws.Range("BN67").Value = 0.6507730628984284
$ws.Range("BP67").Value = 0.65333163591107524

Write-Output "Applied 136 cell updates"
